# Adding 2 search test cases to the "Test Cases" sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Shared-string insertion order matters for an exact diff: descriptions (C) first,
# then Jira ids (B), then TCIDs (A) -- matches how the original sheet was authored.
$ws.Cells.Item(88, 3).Value = "Verify that record view page of a patent gets displayed when user clicks on article title in ALL search results page"
$ws.Cells.Item(89, 3).Value = "Verify that record view page of a patent gets displayed when user clicks a patent title in PATENTS search results page"
$ws.Cells.Item(88, 2).Value = "OPQA-567"
$ws.Cells.Item(89, 2).Value = "OPQA-573"
$ws.Cells.Item(88, 1).Value = "TestCase_B87"
$ws.Cells.Item(89, 1).Value = "TestCase_B88"

$ws.Cells.Item(88, 4).Value = "Y"
$ws.Cells.Item(89, 4).Value = "Y"

# Match the styling used by the rows above (borders / fill / wrap text) for the two new rows
$ws.Range("A87:E87").Copy()
$ws.Range("A88:E89").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Move the active selection/view to mirror the end of the updated list
$ws.Range("C87").Select()
